$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing")

foreach ($r in 36, 81, 122, 167) {
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 1
    $ws.Cells.Item($r, 6).Value = 1
    $ws.Cells.Item($r, 7).Value = 1
    $ws.Cells.Item($r, 8).Value = 1
}

$ws.Range("D36").Select()
